$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '27.656.16'
Set-TextCell 'E2' '  +3.10%  '

# Row 3
Set-TextCell 'D3' '1.853.80'

# Row 4
Set-TextCell 'D4' '1.034'
Set-TextCell 'E4' '  +3.01%  '

# Row 5
Set-TextCell 'D5' '321.61'
Set-TextCell 'E5' '  +4.02%  '

# Row 6
Set-TextCell 'E6' '  +2.70%  '

# Row 7
Set-TextCell 'D7' '0.4380'
Set-TextCell 'E7' '  +1.34%  '

# Row 8
Set-TextCell 'D8' '0.3759'
Set-TextCell 'E8' '  +1.13%  '

# Row 9
Set-TextCell 'E9' '  +2.28%  '

# Row 10
Set-TextCell 'E10' '  +0.97%  '

# Row 11
Set-TextCell 'E11' '  +3.01%  '

# Row 12
Set-TextCell 'D12' '1.865.36'
Set-TextCell 'E12' '  -6.39%  '

# Row 13
Set-TextCell 'E13' '  +2.92%  '

# Row 14
Set-TextCell 'D14' '6.711'
Set-TextCell 'E14' '  +0.61%  '

# Row 15
Set-TextCell 'D15' '0.07158'
Set-TextCell 'E15' '  +3.43%  '

# Row 16
Set-TextCell 'D16' '82.99'
Set-TextCell 'E16' '  +3.04%  '

# Row 17
Set-TextCell 'E17' '  +2.75%  '

# Row 18
Set-TextCell 'D18' '0.000009053'
Set-TextCell 'E18' '  +1.38%  '

# Row 20
Set-TextCell 'D20' '15.45'
Set-TextCell 'E20' '  +1.41%  '

# Row 21
Set-TextCell 'D21' '27.662.02'
Set-TextCell 'E21' '  +2.97%  '

# Row 22
Set-TextCell 'D22' '5.272'
Set-TextCell 'E22' '  +1.12%  '

# Row 23
Set-TextCell 'D23' '11.26'
Set-TextCell 'E23' '  +0.82%  '

# Row 24
Set-TextCell 'D24' '2.076.45'
Set-TextCell 'E24' '  -6.44%  '

# Row 25
Set-TextCell 'D25' '157.64'

# Row 26
Set-TextCell 'D26' '1.944'
Set-TextCell 'E26' '  +3.70%  '

# Row 27
Set-TextCell 'E27' '  +2.41%  '

# Row 28
Set-TextCell 'D28' '5.303'
Set-TextCell 'E28' '  +1.60%  '

# Row 29
Set-TextCell 'D29' '1.941'
Set-TextCell 'E29' '  +1.78%  '

# Row 30
Set-TextCell 'D30' '116.55'
Set-TextCell 'E30' '  +0.97%  '

# Row 31
Set-TextCell 'D31' '0.09086'
Set-TextCell 'E31' '  +1.58%  '

# Row 32
Set-TextCell 'E32' '  +3.11%  '

# Row 33
Set-TextCell 'D33' '0.7690'
Set-TextCell 'E33' '  +1.40%  '

# Row 34
Set-TextCell 'D34' '4.519'
Set-TextCell 'E34' '  +1.67%  '

# Row 35
Set-TextCell 'E35' '  +3.46%  '

# Row 36
Set-TextCell 'E36' '  +2.49%  '

# Row 37
Set-TextCell 'D37' '1.158'
Set-TextCell 'E37' '  +1.77%  '

# Row 38
Set-TextCell 'D38' '0.01980'
Set-TextCell 'E38' '  +2.68%  '

# Row 39
Set-TextCell 'D39' '0.05284'
Set-TextCell 'E39' '  +0.98%  '

# Row 40
Set-TextCell 'B40' 'MXToken'
Set-TextCell 'C40' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D40' '2.826'
Set-TextCell 'E40' '  +6.25%  '

# Row 41
Set-TextCell 'B41' 'TheSandbox'
Set-TextCell 'C41' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D41' '0.5181'
Set-TextCell 'E41' '  +1.76%  '

# Row 42
Set-TextCell 'E42' '  +1.42%  '

# Row 43
Set-TextCell 'D43' '6.734'
Set-TextCell 'E43' '  +2.66%  '

# Row 44
Set-TextCell 'D44' '8.593'
Set-TextCell 'E44' '  +3.62%  '

# Row 45
Set-TextCell 'E45' '  +2.22%  '

# Row 46
Set-TextCell 'D46' '10.61'
Set-TextCell 'E46' '  +1.95%  '

# Row 47
Set-TextCell 'D47' '1.722'
Set-TextCell 'E47' '  +3.79%  '

# Row 48
Set-TextCell 'D48' '0.4661'
Set-TextCell 'E48' '  +2.18%  '

# Row 49
Set-TextCell 'D49' '0.06399'
Set-TextCell 'E49' '  +2.03%  '

# Row 50
Set-TextCell 'E50' '  +4.20%  '

# Row 51
Set-TextCell 'D51' '39.61'
Set-TextCell 'E51' '  +5.58%  '
